$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 675.8333
$ws.Range("J17").Value = 675.8333
$ws.Range("L17").Value = 2027.4999
$ws.Range("N17").Value = -2363.4999

# Row 40
$ws.Range("H40").Value = 1440.8148
$ws.Range("I40").Value = 3725.5
$ws.Range("J40").Value = 1043.4783
$ws.Range("K40").Value = 3725.5
$ws.Range("L40").Value = 1043.4783
$ws.Range("M40").Value = -3550.5
$ws.Range("N40").Value = -1393.4783

# Row 69
$ws.Range("H69").Value = 2866.818
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2866.818
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 8600.454000000002
$ws.Range("N69").Value = -10348.454
$ws.Range("M69").ClearContents()

# Row 72
$ws.Range("H72").Value = 2866.818
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2866.818
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 25801.362
$ws.Range("N72").Value = -34537.362
$ws.Range("M72").ClearContents()

# Row 98
$ws.Range("H98").Value = 1130.4736
$ws.Range("I98").Value = 1130.4736
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1130.4736
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 367.5264
$ws.Range("N98").ClearContents()

# Row 106
$ws.Range("H106").Value = 6064495.5
$ws.Range("I106").Value = 4278.9
$ws.Range("K106").Value = 4278.9
$ws.Range("M106").Value = -3647.9

# Row 122
$ws.Range("H122").Value = 1130.4736
$ws.Range("I122").Value = 1130.4736
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3391.4208
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -941.4207999999999
$ws.Range("N122").ClearContents()

# Row 135
$ws.Range("H135").Value = 21806.541
$ws.Range("I135").Value = 26846.105
$ws.Range("J135").Value = 2656.2
$ws.Range("K135").Value = 241614.945
$ws.Range("L135").Value = 23905.8
$ws.Range("M135").Value = -239079.945
$ws.Range("N135").Value = -28975.8

# Row 137
$ws.Range("H137").Value = 4351390.5
$ws.Range("I137").Value = 5885256.5
$ws.Range("J137").Value = 5437.3335
$ws.Range("K137").Value = 17655769.5
$ws.Range("L137").Value = 16312.0005
$ws.Range("M137").Value = -17653219.5
$ws.Range("N137").Value = -21412.0005

# Row 138
$ws.Range("H138").Value = 3929882.2
$ws.Range("I138").Value = 670518.4
$ws.Range("J138").Value = 5559564
$ws.Range("K138").Value = 2011555.2
$ws.Range("L138").Value = 16678692
$ws.Range("M138").Value = -2006415.2
$ws.Range("N138").Value = -16688972

# Row 141
$ws.Range("H141").Value = 721.3333
$ws.Range("I141").Value = 755.2857
$ws.Range("J141").Value = 602.5
$ws.Range("K141").Value = 2265.8571
$ws.Range("L141").Value = 1807.5
$ws.Range("M141").Value = 2914.1429
$ws.Range("N141").Value = -12167.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1833.3125
$ws.Range("I2").Value = 1594.4166
$ws.Range("K2").Value = 1594.4166
$ws.Range("M2").Value = -1481.4166

# Row 102
$ws.Range("H102").Value = 10206153
$ws.Range("I102").Value = 15874816
$ws.Range("J102").Value = 2560
$ws.Range("K102").Value = 15874816
$ws.Range("L102").Value = 2560
$ws.Range("M102").Value = -15873194
$ws.Range("N102").Value = -5804

# Row 116
$ws.Range("H116").Value = 1833.3125
$ws.Range("I116").Value = 1594.4166
$ws.Range("K116").Value = 1594.4166
$ws.Range("M116").Value = 699.5834

# Row 132
$ws.Range("H132").Value = 71190.07000000001
$ws.Range("I132").Value = 48624
$ws.Range("J132").Value = 123844.22
$ws.Range("K132").Value = 145872
$ws.Range("L132").Value = 371532.66
$ws.Range("M132").Value = -143342
$ws.Range("N132").Value = -376592.66

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1833.3125
$ws.Range("I3").Value = 1594.4166
$ws.Range("K3").Value = 1594.4166
$ws.Range("M3").Value = -1480.4166

# Row 20
$ws.Range("H20").Value = 977.8125
$ws.Range("I20").Value = 661.8
$ws.Range("K20").Value = 661.8
$ws.Range("M20").Value = -414.8

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 43480332
$ws.Range("J58").Value = 2283
$ws.Range("L58").Value = 2283
$ws.Range("N58").Value = -2689

# Row 132
$ws.Range("H132").Value = 12763.453
$ws.Range("I132").Value = 982.16925
$ws.Range("J132").Value = 49229.332
$ws.Range("K132").Value = 2946.50775
$ws.Range("L132").Value = 147687.996
$ws.Range("M132").Value = -416.5077500000002
$ws.Range("N132").Value = -152747.996

# Row 134
$ws.Range("H134").Value = 26935.244
$ws.Range("I134").Value = 2263
$ws.Range("K134").Value = 6789
$ws.Range("M134").Value = -4254

# Row 136
$ws.Range("H136").Value = 43480332
$ws.Range("J136").Value = 2283
$ws.Range("L136").Value = 6849
$ws.Range("N136").Value = -11949

$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Range("H24").Value = 1500
$ws.Range("J24").Value = 1500
$ws.Range("L24").Value = 4500
$ws.Range("N24").Value = -4960

# Row 69
$ws.Range("H69").Value = 875
$ws.Range("J69").Value = 1000
$ws.Range("L69").Value = 3000
$ws.Range("N69").Value = -4622

# Row 72
$ws.Range("H72").Value = 875
$ws.Range("J72").Value = 1000
$ws.Range("L72").Value = 9000
$ws.Range("N72").Value = -17112

# Row 131
$ws.Range("H131").Value = 10870640
$ws.Range("J131").Value = 1123.5581
$ws.Range("L131").Value = 3370.6743
$ws.Range("N131").Value = -13450.6743

# Row 137
$ws.Range("H137").Value = 31998.5
$ws.Range("I137").Value = 1671.6666
$ws.Range("J137").Value = 44995.715
$ws.Range("K137").Value = 5014.9998
$ws.Range("L137").Value = 134987.145
$ws.Range("M137").Value = 85.0002000000004
$ws.Range("N137").Value = -145187.145

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 25222.688
$ws.Range("I70").Value = 42488.08
$ws.Range("J70").Value = 4818.136
$ws.Range("K70").Value = 42488.08
$ws.Range("L70").Value = 4818.136
$ws.Range("M70").Value = -42218.08
$ws.Range("N70").Value = -5358.136

# Row 73
$ws.Range("H73").Value = 25222.688
$ws.Range("I73").Value = 42488.08
$ws.Range("J73").Value = 4818.136
$ws.Range("K73").Value = 42488.08
$ws.Range("L73").Value = 4818.136
$ws.Range("M73").Value = -41552.08
$ws.Range("N73").Value = -6690.136

# Row 122
$ws.Range("H122").Value = 2586.568
$ws.Range("I122").Value = 2302.7354
$ws.Range("J122").Value = 3551.6
$ws.Range("K122").Value = 6908.206200000001
$ws.Range("L122").Value = 10654.8
$ws.Range("M122").Value = -4458.206200000001
$ws.Range("N122").Value = -15554.8

# Row 132
$ws.Range("H132").Value = 73466.5
$ws.Range("I132").Value = 45263.043
$ws.Range("K132").Value = 135789.129
$ws.Range("M132").Value = -133259.129

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3166.6667
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3250
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3250
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -3522

# Row 122
$ws.Range("H122").Value = 3637.75
$ws.Range("I122").Value = 3064.8572
$ws.Range("J122").Value = 4083.3333
$ws.Range("K122").Value = 9194.571599999999
$ws.Range("L122").Value = 12249.9999
$ws.Range("M122").Value = -6744.571599999999
$ws.Range("N122").Value = -17149.9999

# Row 132
$ws.Range("H132").Value = 11090.589
$ws.Range("I132").Value = 13263.409
$ws.Range("J132").Value = 3123.5833
$ws.Range("K132").Value = 39790.227
$ws.Range("L132").Value = 9370.749899999999
$ws.Range("M132").Value = -37260.227
$ws.Range("N132").Value = -14430.7499

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 142860420
$ws.Range("I62").Value = 333334620
$ws.Range("J62").Value = 4750.75
$ws.Range("K62").Value = 333334620
$ws.Range("L62").Value = 4750.75
$ws.Range("M62").Value = -333333996
$ws.Range("N62").Value = -5998.75

# Row 65
$ws.Range("H65").Value = 142860420
$ws.Range("I65").Value = 333334620
$ws.Range("J65").Value = 4750.75
$ws.Range("K65").Value = 1666673100
$ws.Range("L65").Value = 23753.75
$ws.Range("M65").Value = -1666669980
$ws.Range("N65").Value = -29993.75

# Row 126
$ws.Range("H126").Value = 1216.6364
$ws.Range("I126").Value = 1069
$ws.Range("J126").Value = 1475
$ws.Range("K126").Value = 3207
$ws.Range("L126").Value = 4425
$ws.Range("M126").Value = -737
$ws.Range("N126").Value = -9365
